$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset grew by two new price records for "Vega Monumental Concepcion -
# Pimiento". In the canonical row order (sorted by date descending) the two
# new records land right before the former row 594, pushing every row from
# 594..648 down by two (to 596..650).
$ws.Rows.Item(594).Insert()
$ws.Rows.Item(594).Insert()

# New row 594: Cuatro cascos verde, Primera, Provincia de Limari
$ws.Range("A594").Value = 11
$ws.Range("B594").Value = "Vega Monumental Concepción"
$ws.Range("C594").Value = "Bíobío"
$ws.Range("D594").Value = 45223
$ws.Range("E594").Value = 8
$ws.Range("F594").Value = 100112002
$ws.Range("G594").Value = "Pimiento"
$ws.Range("H594").Value = "Cuatro cascos verde"
$ws.Range("I594").Value = "Primera"
$ws.Range("J594").Value = 120
$ws.Range("K594").Value = 30000
$ws.Range("L594").Value = 30000
$ws.Range("M594").Value = 30000
$ws.Range("N594").Value = "$/caja 18 kilos"
$ws.Range("O594").Value = "Provincia de Limarí"
$ws.Range("P594").Value = 1667
$ws.Range("Q594").Value = 18
$ws.Range("R594").Value = "Hortaliza"

# New row 595: Zafiro rojo, Primera, Región de Arica y Parinacota
$ws.Range("A595").Value = 11
$ws.Range("B595").Value = "Vega Monumental Concepción"
$ws.Range("C595").Value = "Bíobío"
$ws.Range("D595").Value = 45223
$ws.Range("E595").Value = 8
$ws.Range("F595").Value = 100112002
$ws.Range("G595").Value = "Pimiento"
$ws.Range("H595").Value = "Zafiro rojo"
$ws.Range("I595").Value = "Primera"
$ws.Range("J595").Value = 100
$ws.Range("K595").Value = 30000
$ws.Range("L595").Value = 30000
$ws.Range("M595").Value = 30000
$ws.Range("N595").Value = "$/caja 15 kilos"
$ws.Range("O595").Value = "Región de Arica y Parinacota"
$ws.Range("P595").Value = 2000
$ws.Range("Q595").Value = 15
$ws.Range("R595").Value = "Hortaliza"
